# Switch license from BY-NC to BY-SA
# (also tidies up the title-slide heading into a single run)

$p = $ppt.ActivePresentation

# --- Slide 1: title slide -------------------------------------------------
# "Blue Waters Petascale" + " Semester Curriculum v1.0" -> one run
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(1)
$tr1 = $sh1.TextFrame.TextRange

$full1 = $tr1.Text
$combined = "Blue Waters Petascale Semester Curriculum v1.0"
$idx1 = $full1.IndexOf("Blue Waters Petascale")
$titleRange = $tr1.Characters($idx1 + 1, $combined.Length)
$titleRange.Text = $combined

# --- Slide 2: license / closing slide -------------------------------------
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(1)
$tr2 = $sh2.TextFrame.TextRange

# Nudge the title placeholder's horizontal position (566059 -> 566057 EMU)
$sh2.Left = 44.57141732283465

# "CC BY-NC 4.0. ..." -> "CC BY-SA 4.0. ..."
$full2 = $tr2.Text
$byncIdx = $full2.IndexOf("BY-NC ")
$byncRange = $tr2.Characters($byncIdx + 1, "BY-NC ".Length)
$byncRange.Text = "BY-SA "

# "https://creativecommons.org/licenses/by-nc/4.0" -> ".../by-sa/4.0",
# then split the run so the "https://" prefix and the rest are separate runs
$full2 = $tr2.Text
$oldUrl = "https://creativecommons.org/licenses/by-nc/4.0"
$newUrl = "https://creativecommons.org/licenses/by-sa/4.0"
$urlIdx = $full2.IndexOf($oldUrl)
$urlRange = $tr2.Characters($urlIdx + 1, $oldUrl.Length)
$urlRange.Text = $newUrl

$full2 = $tr2.Text
$newUrlIdx = $full2.IndexOf($newUrl)
$prefixRange = $tr2.Characters($newUrlIdx + 1, "https://".Length)
$prefixRange.Text = "https://"
